$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.292789816856384
$ws.Range("B1").Value = 2.204940795898438
$ws.Range("D1").Value = 1.385526776313782
$ws.Range("E1").Value = 0.8221459984779358
